$d = $word.ActiveDocument

$pairs = @(
    @("363÷7=51, 6", "981÷6=163, 3"),
    @("537÷6=89, 3", "347÷7=49, 4"),
    @("226÷8=28, 2", "332÷7=47, 3"),
    @("368÷6=61, 2", "484÷6=80, 4"),
    @("118÷2=59, 0", "235÷9=26, 1"),
    @("829÷8=103, 5", "186÷3=62, 0"),
    @("146÷2=73, 0", "401÷8=50, 1"),
    @("868÷3=289, 1", "340÷7=48, 4"),
    @("145÷2=72, 1", "116÷9=12, 8"),
    @("138÷3=46, 0", "194÷2=97, 0"),
    @("873÷8=109, 1", "165÷6=27, 3"),
    @("116÷5=23, 1", "953÷3=317, 2"),
    @("245÷4=61, 1", "353÷2=176, 1"),
    @("584÷5=116, 4", "219÷9=24, 3"),
    @("109÷8=13, 5", "124÷6=20, 4"),
    @("147÷3=49, 0", "337÷5=67, 2"),
    @("506÷9=56, 2", "637÷3=212, 1"),
    @("704÷8=88, 0", "959÷2=479, 1"),
    @("874÷3=291, 1", "688÷3=229, 1"),
    @("745÷2=372, 1", "440÷9=48, 8"),
    @("351÷6=58, 3", "290÷7=41, 3"),
    @("607÷4=151, 3", "514÷7=73, 3"),
    @("661÷4=165, 1", "937÷3=312, 1"),
    @("853÷5=170, 3", "940÷4=235, 0"),
    @("594÷3=198, 0", "628÷5=125, 3")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
